$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset (rows 329-364) gains one more weekly reading (Primera / Segunda)
# at the front of that block. Shift the existing rows 329-364 down by two rows
# (native Excel row-insert semantics: formatting + all columns move together),
# then populate the two freshly inserted rows with the new week's data.
$ws.Range("A329:R330").EntireRow.Insert()

# New row 329 - "Primera"
$ws.Cells.Item(329, 1).Value = 3
$ws.Cells.Item(329, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 44449
$ws.Cells.Item(329, 5).Value = 5
$ws.Cells.Item(329, 6).Value = 100114014
$ws.Cells.Item(329, 7).Value = "Betarraga"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Primera"
$ws.Cells.Item(329, 10).Value = 3100
$ws.Cells.Item(329, 11).Value = 500
$ws.Cells.Item(329, 12).Value = 550
$ws.Cells.Item(329, 13).Value = 524
$ws.Cells.Item(329, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(329, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(329, 16).Value = 131
$ws.Cells.Item(329, 17).Value = 4
$ws.Cells.Item(329, 18).Value = "Hortaliza"

# New row 330 - "Segunda"
$ws.Cells.Item(330, 1).Value = 3
$ws.Cells.Item(330, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(330, 3).Value = "Coquimbo"
$ws.Cells.Item(330, 4).Value = 44449
$ws.Cells.Item(330, 5).Value = 5
$ws.Cells.Item(330, 6).Value = 100114014
$ws.Cells.Item(330, 7).Value = "Betarraga"
$ws.Cells.Item(330, 8).Value = "Sin especificar"
$ws.Cells.Item(330, 9).Value = "Segunda"
$ws.Cells.Item(330, 10).Value = 1200
$ws.Cells.Item(330, 11).Value = 400
$ws.Cells.Item(330, 12).Value = 400
$ws.Cells.Item(330, 13).Value = 400
$ws.Cells.Item(330, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(330, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(330, 16).Value = 100
$ws.Cells.Item(330, 17).Value = 4
$ws.Cells.Item(330, 18).Value = "Hortaliza"
